$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text content looks numeric need a temporary "Text" number
# format applied before the write, otherwise the COM layer auto-coerces
# them into numeric cells (losing the shared-string type / leading zeros).
$ws.Range("B7").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("G7").NumberFormat = "@"

$ws.Range("A7").Value = "Вильям Иванович Шекспир"
$ws.Range("B7").Value = "-666"
$ws.Range("C7").Value = "20-05-1456"
$ws.Range("D7").Value = "29654"
$ws.Range("E7").Value = "Гамлет"
$ws.Range("F7").Value = "1шт"
$ws.Range("G7").Value = "006312"
$ws.Range("H7").Value = "Долго не живет"
$ws.Range("I7").Value = "-"

# Restore default (General) styling on those cells so the written row
# matches the rest of the sheet (no explicit style index).
$ws.Range("B7").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("G7").ClearFormats()
